$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 19.51877066666666
$ws.Range("N2").Value = 58.556312
$ws.Range("O2").Value = 0.9213325069349741
$ws.Range("P2").Value = 0.9213325069349743
$ws.Range("Q2").Value = 0.8503157253226665
$ws.Range("R2").Value = 7.652841527904
$ws.Range("S2").Value = 0.9213325069349741
$ws.Range("T2").Value = 0.9213325069349743

# Row 3
$ws.Range("O3").Value = 0.01004540680524951
$ws.Range("P3").Value = 0.01004540680524951
$ws.Range("S3").Value = 0.01004540680524951
$ws.Range("T3").Value = 0.01004540680524951

# Row 4
$ws.Range("N4").Value = 4.361353
$ws.Range("O4").Value = 0.06862208625977624
$ws.Range("P4").Value = 0.06862208625977625
$ws.Range("R4").Value = 0.5699939462760001
$ws.Range("S4").Value = 0.06862208625977624
$ws.Range("T4").Value = 0.06862208625977625
